$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.315.43"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "1.552.70"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'209.89"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'23.75"
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "1.774.10"
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("D13").Value = "1.548.06"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").Value = "28.289.01"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.509"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'3.61"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "'60.55"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").Value = "'227.67"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'7.32"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "0.0₃0674"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'3.91"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "'8.89"
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("E24").Value = "  -4.18%  "
$ws.Range("D25").Value = "'151.67"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").Value = "'14.75"
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("E30").Value = "  -3.64%  "
$ws.Range("E31").Value = "  -4.83%  "
$ws.Range("D32").Value = "'3.16"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("D33").Value = "1.386.66"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "'3.02"
$ws.Range("E34").Value = "  -2.96%  "
$ws.Range("D35").Value = "'1.07"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "'1.48"
$ws.Range("E36").Value = "  -4.22%  "
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("D41").Value = "'0.510"
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("D46").Value = "'61.88"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").Value = "1.687.23"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("E48").Value = "  -6.17%  "
$ws.Range("D49").Value = "'85.68"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").Value = "'42.19"
$ws.Range("E50").Value = "  +6.32%  "
$ws.Range("E51").Value = "  +6.77%  "
